# Generate Report for Handoff
# Rotates the localization-status report's generated GUID/hash and
# refreshes the handoff timestamps, mirroring a fresh CI run.

$wb = $excel.ActiveWorkbook

$oldGuid = "38fcd06a-0aab-48cc-9898-9f7d7e84acc5"
$newGuid = "a66e5554-e7d1-42ca-acc7-7540ff62b6d0"
$oldHash = "0e8c66277725376d93d47c483d0e2f2c53dde5f5"
$newHash = "d4f3382de2e32a769bd1413e389a30c064bc2241"

$newMdName = "$newGuid.md"
$newMdPath = "e2e\$newGuid.md"

$overviewDate = "2016-09-03 09:01:15"
$zhHandoffFile = "$newGuid.$newHash.zh-cn.xlf"
$zhHandoffDate = "2016-09-03 09:01:10"
$deHandoffFile = "$newGuid.$newHash.de-de.xlf"
$deHandoffDate = $overviewDate

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Cells.Item(2, 1).Value = $newMdName      # A2 File Name
$wsOverview.Cells.Item(2, 2).Value = $newMdPath      # B2 Path And Name
$wsOverview.Cells.Item(2, 7).Value = $overviewDate   # G2 Latest HO Xliff Generate Date

foreach ($h in $wsOverview.Hyperlinks) {
    $h.TextToDisplay = $newMdPath
}

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Cells.Item(2, 1).Value = $newMdName        # A2 Source File Name
$wsZh.Cells.Item(2, 7).Value = $zhHandoffFile     # G2 Latest Handoff File
$wsZh.Cells.Item(2, 8).Value = $zhHandoffDate     # H2 Latest Handoff Datetime

foreach ($h in $wsZh.Hyperlinks) {
    $h.TextToDisplay = $newMdName
}

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Cells.Item(2, 1).Value = $newMdName        # A2 Source File Name
$wsDe.Cells.Item(2, 7).Value = $deHandoffFile     # G2 Latest Handoff File
$wsDe.Cells.Item(2, 8).Value = $deHandoffDate     # H2 Latest Handoff Datetime (shared w/ Overview date)

foreach ($h in $wsDe.Hyperlinks) {
    $h.TextToDisplay = $newMdName
}
